# Restyle the deck's theme: swap the "Integral" theme colours that are
# currently applied to the slide master for the stock "Office Theme"
# colour palette (ppt/theme/theme1.xml in the package).
#
# Helper: turn an "RRGGBB" hex string into the packed BGR decimal value
# that the PowerPoint OLE `RGB` color properties expect.
function HexToOleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# The 12 slots of the theme colour scheme, in PowerPoint's fixed order:
# Dark1, Light1, Dark2, Light2, Accent1-6, Hyperlink, FollowedHyperlink.
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $colorScheme.Item($i + 1).RGB = HexToOleColor $officeThemeColors[$i]
}
